$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are stored as text so that numeric-looking
# strings such as "69.405.21" or "1.00" are preserved exactly as text,
# matching the source data which never contains genuine numeric cells here.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '69.405.21'
$ws.Range("E2").Value = '  +2.21%  '
$ws.Range("D3").Value = '3.403.02'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '588.04'
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("D6").Value = '181.56'
$ws.Range("E6").Value = '  +3.70%  '
$ws.Range("D7").Value = '0.600'
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("D9").Value = '0.203'
$ws.Range("E9").Value = '  +11.54%  '
$ws.Range("E10").Value = '  +2.37%  '
$ws.Range("D11").Value = '48.54'
$ws.Range("E11").Value = '  +2.65%  '
$ws.Range("E12").Value = '  +5.41%  '
$ws.Range("D13").Value = '686.16'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").Value = '3.952.37'
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").Value = '69.454.32'
$ws.Range("E16").Value = '  +2.30%  '
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '3.402.07'
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("D19").Value = '17.75'
$ws.Range("E19").Value = '  +1.52%  '
$ws.Range("D20").Value = '11.35'
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = '0.912'
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("D22").Value = '17.36'
$ws.Range("E22").Value = '  +2.52%  '
$ws.Range("D23").Value = '5.37'
$ws.Range("E23").Value = '  -0.96%  '
$ws.Range("D24").Value = '103.39'
$ws.Range("E24").Value = '  +2.10%  '
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("E26").Value = '  +1.65%  '
$ws.Range("D27").Value = '9.72'
$ws.Range("E27").Value = '  +2.79%  '
$ws.Range("D28").Value = '34.01'
$ws.Range("D29").Value = '8.85'
$ws.Range("E29").Value = '  +3.67%  '
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").Value = '11.17'
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("E32").Value = '  -2.84%  '
$ws.Range("D33").Value = '3.64'
$ws.Range("E33").Value = '  +11.04%  '
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").Value = '58.62'
$ws.Range("E35").Value = '  +3.48%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("D37").Value = '3.660.03'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("E38").Value = '  +5.75%  '
$ws.Range("D39").Value = '35.99'
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").Value = '  +8.39%  '
$ws.Range("D41").Value = '3.27'
$ws.Range("E41").Value = '  +4.19%  '
$ws.Range("D42").Value = '2.68'
$ws.Range("D43").Value = '0.0429'
$ws.Range("E43").Value = '  +5.69%  '
$ws.Range("D44").Value = '0.339'
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.35'
$ws.Range("E45").Value = '  +1.51%  '
$ws.Range("B46").Value = 'ThetaToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D46").Value = '2.68'
$ws.Range("E46").Value = '  +2.42%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '0.130'
$ws.Range("E47").Value = '  +1.00%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '1.39'
$ws.Range("E48").Value = '  +5.05%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '131.06'
$ws.Range("E50").Value = '  +0.36%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").Value = '2.74'
$ws.Range("E51").Value = '  +1.90%  '
